$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 34
$ws.Range("C2").Value = 54
$ws.Range("C3").Value = 34

$ws.Range("C4").Select()
